# Automatische test-sync: 2025-07-22 12:52:50
# Adds Testmail #19 to the "Logs" sheet (row 18) and syncs the category
# counts on the "Dashboard" sheet accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Logs" sheet: append the new mail log entry as row 18.
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A18").Value = "Heeft u informatie over zakelijke kortingen voor wederverkopers?"
$logs.Range("B18").Value = "mailmind.test@zohomail.eu"
$logs.Range("C18").Value = "Testmail #19: Heeft u informatie over zakelijke kortingen voor wederverkopers?"
$logs.Range("D18").Value = "Productinformatie"
$logs.Range("E18").Value = "Beste klant,`nBedankt voor uw interesse in onze zakelijke kortingen voor wederverkopers. Voor meer informatie over onze zakelijke kortingen en de voorwaarden kunt u contact opnemen met onze verkoopafdeling via sales@bedrijfsnaam.nl. Zij kunnen u voorzien van alle benodigde informatie en u verder helpen met uw aanvraag.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Range("F18").Value = "2025-07-22 12:52:47"
$logs.Range("G18").Value = "Ja"
$logs.Range("H18").Value = "Nee"
$logs.Range("I18").Value = "Ja"
$logs.Range("J18").Value = "Ja"

# Extend the conditional-formatting ranges from row 17 to row 18 so the new
# row picks up the same category / Ja-Nee highlighting.
$logs.Range("D2:D17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D18"))
$logs.Range("G2:G17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G18"))
$logs.Range("H2:H17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H18"))
$logs.Range("I2:I17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I18"))
$logs.Range("J2:J17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J18"))

# ---------------------------------------------------------------------------
# 2. "Dashboard" sheet: the new mail bumps "Productinformatie" from 4 to 5,
#    tying it with "Retour / Terugbetaling" (still 5). The two rows swap
#    places so the list stays ordered by descending count.
# ---------------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A2").Value = "Productinformatie"
$dashboard.Range("B2").Value = 5
$dashboard.Range("A3").Value = "Retour / Terugbetaling"
$dashboard.Range("B3").Value = 5
